# The bulleted list under "Створення замовлення" originally read:
#   1. Додавання запису у таблицю Orders
#   2. Додавання запису у таблицю DeliveryInfos
#   3. Додавання записів у таблицю FurnitureOrderRows
#   4. Додавання записів у таблицю AdditionalDetailsOrdered
#   5. Прорахування ціни кожного рядку замовлення
#   6. Прорахування загальної ціни замовлення
#   7. Прорахування заробленої суми із замволення та занесення її в таблицю Profits
#
# The edit moves the three "calculation" bullets (5,6,7) to the front of the
# list (ahead of the "add a record" bullets) and fixes a typo
# ("замволення" -> "замовлення") in what is now the 3rd bullet. Net effect
# on list order: [1,2,3,4,5,6,7] -> [5,6,7,1,2,3,4].
#
# We reproduce this with whole-paragraph Cut/Paste so each paragraph keeps
# its original run structure (incl. the run split between the descriptive
# text and the bare table-name run) instead of retyping text, which would
# risk merging/losing the run boundaries.

$d = $word.ActiveDocument

# The list starts right after the "Створення замовлення" title paragraph
# (paragraph 1), so bullets 1..7 are currently Word paragraphs 2..8.
# Repeatedly cut the paragraph currently in slot 8 (the next "calculation"
# bullet working from the back) and paste it right after paragraph 1.
# Doing this three times in a row yields, in order: bullet7, then bullet6,
# then bullet5 each inserted right after the title -- which leaves the
# list in the order bullet5, bullet6, bullet7, bullet1, bullet2, bullet3,
# bullet4, exactly the target order.
for ($i = 0; $i -lt 3; $i++) {
    $cutPara = $d.Paragraphs.Item(8)
    $moveRange = $d.Range($cutPara.Range.Start, $cutPara.Range.End)
    $moveRange.Cut() | Out-Null

    $titleEnd = $d.Paragraphs.Item(1).Range.End
    $insertRange = $d.Range($titleEnd, $titleEnd)
    $insertRange.Paste() | Out-Null
}

# Fix the typo in the "Profits" bullet (now the 3rd bullet in the list):
# "замволення" -> "замовлення".
$d.Content.Find.Execute("із замволення та занесення", $false, $false, $false, $false, $false, $true, 1, $false, "із замовлення та занесення", 2) | Out-Null
